# Adds a new worksheet table of Common / Scientific fish names on "Sheet2",
# mirroring the indicator species already listed on "Sheet1".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- header row -------------------------------------------------------------
$ws2.Range("B1").Value = "Common"
$ws2.Range("C1").Value = "Scientific"

# --- data rows (Common name / Scientific name) ------------------------------
$ws2.Range("B2").Value  = "Bluehead Chub"
$ws2.Range("C2").Value  = "Nocomis_leptocephalus"

$ws2.Range("B3").Value  = "Fantail Darter/Carolina Fantail Darter"
$ws2.Range("C3").Value  = "Etheostoma_flabellare"

$ws2.Range("B4").Value  = "Pirate Perch"
$ws2.Range("C4").Value  = "Aphredoderus_sayanus"

$ws2.Range("B5").Value  = "Margined Madtom"
$ws2.Range("C5").Value  = "Noturus_insignis"

$ws2.Range("B6").Value  = "Notchlip Redhorse"
$ws2.Range("C6").Value  = "Moxostoma_collapsum"

$ws2.Range("B7").Value  = "Redlip Shiner (where native)/Greenhead Shiner/Piedmont Shiner"

$ws2.Range("B8").Value  = "Whitemouth Shiner/Swallowtail Shiner (where native)"

$ws2.Range("B9").Value  = "Chainback Darter/Piedmont Darter"
$ws2.Range("C9").Value  = "Percina_nevisense"

$ws2.Range("B10").Value = "Redbreast Sunfish"
$ws2.Range("C10").Value = "Lepomis_auritus"

$ws2.Range("B11").Value = "Highfin Shiner"
$ws2.Range("C11").Value = "Notropis_altipinnis"

$ws2.Range("B12").Value = "Rosyside Dace"
$ws2.Range("C12").Value = "Clinostomus_funduloides"

$ws2.Range("B13").Value = "Tessellated Darter/Johnny Darter"
$ws2.Range("C13").Value = "Etheostoma_olmstedi"

# --- reuse the "species name" cell format already used for this data on
#     Sheet1 (Arial 10, vertically centered) instead of re-deriving it -------
$ws1.Range("C13:C24").Copy()
$ws2.Range("B2:B13").PasteSpecial(-4122)   # xlPasteFormats

# --- column widths, matching Sheet1's corresponding columns -----------------
$ws2.Columns.Item(2).ColumnWidth = 55.140625
$ws2.Columns.Item(3).ColumnWidth = 31.7109375

# --- selection / active-sheet bookkeeping ------------------------------------
# Sheet1 was left with C13:C24 selected before switching away from it.
$ws1.Range("C13:C24").Select()

# Sheet2 becomes the active (visible) tab, with C13 selected.
$ws2.Activate()
$ws2.Range("C13").Select()
